$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest refresh
$ws.Range("D2").Value = "26.712.95"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.602.06"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.68"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0620"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.247"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.73"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "1.826.07"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "1.596.01"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.26"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").Value = "26.687.49"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "0.0₃0747"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.26"
$ws.Range("E19").Value = "  +2.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.68"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.40"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.13"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.43"
$ws.Range("E29").Value = "  +1.27%  "
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").Value = "1.298.67"
$ws.Range("E34").Value = "  +2.10%  "
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.610"
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("E38").Value = "  +22.08%  "
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.41"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.20"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.782"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.41"
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("D45").Value = "1.737.08"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.23"
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("E47").Value = "  -2.24%  "
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0518"
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("E51").Value = "  -0.22%  "
